$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 199
$ws.Range("I6").Value = 115.76923
$ws.Range("J6").Value = 379.33334
$ws.Range("K6").Value = 347.30769
$ws.Range("L6").Value = 1138.00002
$ws.Range("M6").Value = -235.30769
$ws.Range("N6").Value = -1362.00002
# Row 9
$ws.Range("H9").Value = 421.36365
$ws.Range("I9").Value = 307.2
$ws.Range("K9").Value = 307.2
$ws.Range("M9").Value = -138.2
# Row 31
$ws.Range("H31").Value = 249.5
$ws.Range("I31").Value = 249.5
$ws.Range("K31").Value = 748.5
$ws.Range("M31").Value = -518.5
# Row 40
$ws.Range("H40").Value = 5533.25
$ws.Range("J40").Value = 5928.4287
$ws.Range("L40").Value = 5928.4287
$ws.Range("N40").Value = -6278.4287
# Row 74
$ws.Range("H74").Value = 9454.091
$ws.Range("J74").Value = 14750
$ws.Range("L74").Value = 14750
$ws.Range("N74").Value = -16622
# Row 77
$ws.Range("H77").Value = 9454.091
$ws.Range("J77").Value = 14750
$ws.Range("L77").Value = 73750
$ws.Range("N77").Value = -83110
# Row 80
$ws.Range("H80").Value = 3684.9167
$ws.Range("I80").Value = 2087.25
$ws.Range("K80").Value = 6261.75
$ws.Range("M80").Value = -5263.75
# Row 83
$ws.Range("H83").Value = 3684.9167
$ws.Range("I83").Value = 2087.25
$ws.Range("K83").Value = 18785.25
$ws.Range("M83").Value = -13793.25
# Row 116
$ws.Range("H116").Value = 10690.5
$ws.Range("I116").Value = 10499.5
$ws.Range("J116").Value = 10738.25
$ws.Range("K116").Value = 10499.5
$ws.Range("L116").Value = 10738.25
$ws.Range("M116").Value = -7057.5
$ws.Range("N116").Value = -17622.25

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2500
$ws.Range("I61").Value = 1500
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1500
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1288
$ws.Range("N61").Value = -3424
# Row 102
$ws.Range("H102").Value = 3318.8
$ws.Range("I102").Value = 3606.5
$ws.Range("J102").Value = 729.5
$ws.Range("K102").Value = 3606.5
$ws.Range("L102").Value = 729.5
$ws.Range("M102").Value = -1984.5
$ws.Range("N102").Value = -3973.5
# Row 110
$ws.Range("H110").Value = 4583.476
$ws.Range("I110").Value = 4562.65
$ws.Range("K110").Value = 4562.65
$ws.Range("M110").Value = -2517.65
# Row 122
$ws.Range("H122").Value = 3887.2144
$ws.Range("I122").Value = 3320.4583
$ws.Range("K122").Value = 9961.374899999999
$ws.Range("M122").Value = -7511.374899999999
# Row 132
$ws.Range("H132").Value = 2876.647
$ws.Range("I132").Value = 2868.1428
$ws.Range("J132").Value = 2916.3333
$ws.Range("K132").Value = 8604.428400000001
$ws.Range("L132").Value = 8748.999899999999
$ws.Range("M132").Value = -6074.428400000001
$ws.Range("N132").Value = -13808.9999
# Row 136
$ws.Range("H136").Value = 2500
$ws.Range("I136").Value = 1500
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 4500
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -1950
$ws.Range("N136").Value = -14100

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 184.75
$ws.Range("I22").Value = 101
$ws.Range("J22").Value = 436
$ws.Range("K22").Value = 101
$ws.Range("L22").Value = 436
$ws.Range("M22").Value = 72
$ws.Range("N22").Value = -782
# Row 94
$ws.Range("H94").Value = 2108.457
$ws.Range("I94").Value = 1804.64
$ws.Range("K94").Value = 1804.64
$ws.Range("M94").Value = -1353.64
# Row 105
$ws.Range("H105").Value = 2229.7144
$ws.Range("I105").Value = 1689.5
$ws.Range("K105").Value = 1689.5
$ws.Range("M105").Value = 57.5
# Row 134
$ws.Range("H134").Value = 3259.2856
$ws.Range("I134").Value = 2953.75
$ws.Range("J134").Value = 3666.6667
$ws.Range("K134").Value = 8861.25
$ws.Range("L134").Value = 11000.0001
$ws.Range("M134").Value = -6326.25
$ws.Range("N134").Value = -16070.0001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
# Row 132
$ws.Range("H132").Value = 3271.2942
$ws.Range("I132").Value = 2601.5454
$ws.Range("K132").Value = 7804.6362
$ws.Range("M132").Value = -5274.6362

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 2814.1428
$ws.Range("I68").Value = 1066.3334
$ws.Range("K68").Value = 3199.0002
$ws.Range("M68").Value = -2388.0002
# Row 71
$ws.Range("H71").Value = 2814.1428
$ws.Range("I71").Value = 1066.3334
$ws.Range("K71").Value = 9597.000599999999
$ws.Range("M71").Value = -5541.000599999999
# Row 132
$ws.Range("H132").Value = 2598.25
$ws.Range("J132").Value = 2200
$ws.Range("L132").Value = 19800
$ws.Range("N132").Value = -24860

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 289529.75
$ws.Range("I70").Value = 377706.34
$ws.Range("K70").Value = 377706.34
$ws.Range("M70").Value = -377436.34
# Row 73
$ws.Range("H73").Value = 289529.75
$ws.Range("I73").Value = 377706.34
$ws.Range("K73").Value = 377706.34
$ws.Range("M73").Value = -376770.34
# Row 122
$ws.Range("H122").Value = 6382.385
$ws.Range("I122").Value = 6627.4287
$ws.Range("K122").Value = 19882.2861
$ws.Range("M122").Value = -17432.2861

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2336.6191
$ws.Range("I132").Value = 1833.1765
$ws.Range("K132").Value = 5499.529500000001
$ws.Range("M132").Value = -2969.529500000001
# Row 136
$ws.Range("H136").Value = 4394.8335
$ws.Range("I136").Value = 4025.0527
$ws.Range("J136").Value = 5800
$ws.Range("K136").Value = 12075.1581
$ws.Range("L136").Value = 17400
$ws.Range("M136").Value = -9525.158100000001
$ws.Range("N136").Value = -22500

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1332.375
$ws.Range("I132").Value = 858.85
$ws.Range("J132").Value = 3700
$ws.Range("K132").Value = 2576.55
$ws.Range("L132").Value = 11100
$ws.Range("M132").Value = -46.55000000000018
$ws.Range("N132").Value = -16160
# Row 136
$ws.Range("H136").Value = 1190.2727
$ws.Range("I136").Value = 934
$ws.Range("J136").Value = 1638.75
$ws.Range("K136").Value = 2802
$ws.Range("L136").Value = 4916.25
$ws.Range("M136").Value = -252
$ws.Range("N136").Value = -10016.25

Write-Output "Applied all Ultros_Profits market-data updates"